$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('E2').Value = '2026-02-13 21:48:40'
$ws.Range('I2').Value = '3.8 mm'
$ws.Range('E3').Value = '2026-02-13 21:48:43'
$ws.Range('I3').Value = '7.3 mm'
$ws.Range('E4').Value = '2026-02-13 21:48:45'
$ws.Range('J4').Value = '993.6 hPa'
$ws.Range('L4').Value = '22.7 km/h - 287º 21:22 TU'
$ws.Range('E5').Value = '2026-02-13 21:48:48'
$ws.Range('G5').Value = '110 cm'
$ws.Range('H5').Value = "'83%"
$ws.Range('I5').Value = '3.2 mm'
$ws.Range('O5').Value = '-2.8 °C'
$ws.Range('E6').Value = '2026-02-13 21:48:50'
$ws.Range('J6').Value = '993.6 hPa'
$ws.Range('E7').Value = '2026-02-13 21:48:53'
$ws.Range('J7').Value = '994.0 hPa'
$ws.Range('L7').Value = '57.2 km/h - 314º 21:26 TU'
$ws.Range('E8').Value = '2026-02-13 21:48:56'
$ws.Range('J8').Value = '993.9 hPa'
$ws.Range('L8').Value = '50.0 km/h - 296º 21:28 TU'
$ws.Range('E9').Value = '2026-02-13 21:48:58'
$ws.Range('E10').Value = '2026-02-13 21:49:00'
$ws.Range('E11').Value = '2026-02-13 21:49:03'
$ws.Range('E12').Value = '2026-02-13 21:49:05'
$ws.Range('O12').Value = '9.5 °C'
$ws.Range('E13').Value = '2026-02-13 21:49:07'
$ws.Range('J13').Value = '996.7 hPa'
$ws.Range('E14').Value = '2026-02-13 21:49:10'
$ws.Range('L14').Value = '39.2 km/h - 319º 21:15 TU'
$ws.Range('O14').Value = '10.5 °C'
$ws.Range('E15').Value = '2026-02-13 21:49:12'
$ws.Range('E16').Value = '2026-02-13 21:49:15'
$ws.Range('H16').Value = "'80%"
$ws.Range('I16').Value = '13.8 mm'
$ws.Range('E17').Value = '2026-02-13 21:49:17'
$ws.Range('H17').Value = "'91%"
$ws.Range('E18').Value = '2026-02-13 21:49:20'
$ws.Range('J18').Value = '993.8 hPa'
$ws.Range('E19').Value = '2026-02-13 21:49:22'
$ws.Range('E20').Value = '2026-02-13 21:49:25'
$ws.Range('I20').Value = '24.2 mm'
$ws.Range('E21').Value = '2026-02-13 21:49:27'
$ws.Range('J21').Value = '996.8 hPa'
$ws.Range('N21').Value = '-0.2 °C 21:11 TU'
$ws.Range('E22').Value = '2026-02-13 21:49:29'
$ws.Range('L22').Value = '49.7 km/h - 337º 21:21 TU'
$ws.Range('E23').Value = '2026-02-13 21:49:32'
$ws.Range('I23').Value = '12.2 mm'
$ws.Range('O23').Value = '-4.0 °C'
$ws.Range('E24').Value = '2026-02-13 21:49:34'
$ws.Range('J24').Value = '994.9 hPa'
$ws.Range('L24').Value = '61.6 km/h - 294º 21:08 TU'
$ws.Range('E25').Value = '2026-02-13 21:49:37'
$ws.Range('I25').Value = '9.6 mm'
$ws.Range('O25').Value = '-2.8 °C'
$ws.Range('E26').Value = '2026-02-13 21:49:39'
$ws.Range('E27').Value = '2026-02-13 21:49:42'
$ws.Range('E28').Value = '2026-02-13 21:49:44'
$ws.Range('J28').Value = '994.1 hPa'
$ws.Range('E29').Value = '2026-02-13 21:49:47'
$ws.Range('E30').Value = '2026-02-13 21:49:49'
$ws.Range('J30').Value = '993.6 hPa'
$ws.Range('E31').Value = '2026-02-13 21:49:52'
$ws.Range('H31').Value = "'75%"
$ws.Range('J31').Value = '992.5 hPa'
$ws.Range('E32').Value = '2026-02-13 21:49:54'
$ws.Range('L32').Value = '58.0 km/h - 310º 21:26 TU'
$ws.Range('O32').Value = '5.0 °C'
$ws.Range('E33').Value = '2026-02-13 21:49:57'
$ws.Range('J33').Value = '995.7 hPa'
$ws.Range('E34').Value = '2026-02-13 21:49:59'
$ws.Range('O34').Value = '-0.7 °C'
$ws.Range('E35').Value = '2026-02-13 21:50:02'
$ws.Range('J35').Value = '995.0 hPa'
$ws.Range('L35').Value = '83.2 km/h - 257º 21:11 TU'
$ws.Range('E36').Value = '2026-02-13 21:50:04'
$ws.Range('J36').Value = '993.7 hPa'
$ws.Range('E37').Value = '2026-02-13 21:50:07'
$ws.Range('J37').Value = '995.6 hPa'
$ws.Range('E38').Value = '2026-02-13 21:50:09'
$ws.Range('H38').Value = "'79%"
$ws.Range('E39').Value = '2026-02-13 21:50:12'
$ws.Range('I39').Value = '19.7 mm'
$ws.Range('O39').Value = '-3.7 °C'
$ws.Range('E40').Value = '2026-02-13 21:50:14'
$ws.Range('J40').Value = '997.3 hPa'
$ws.Range('E41').Value = '2026-02-13 21:50:17'
$ws.Range('H41').Value = "'75%"
$ws.Range('J41').Value = '994.3 hPa'
$ws.Range('E42').Value = '2026-02-13 21:50:19'
$ws.Range('E43').Value = '2026-02-13 21:50:22'
$ws.Range('E44').Value = '2026-02-13 21:50:24'
$ws.Range('I44').Value = '9.9 mm'
$ws.Range('E45').Value = '2026-02-13 21:50:27'
$ws.Range('H45').Value = "'67%"
$ws.Range('I45').Value = '2.5 mm'
$ws.Range('J45').Value = '993.0 hPa'
$ws.Range('E46').Value = '2026-02-13 21:50:29'
$ws.Range('H46').Value = "'88%"
$ws.Range('J46').Value = '995.1 hPa'
$ws.Range('L46').Value = '60.5 km/h - 320º 21:23 TU'
